$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp in A1 (07:04 -> 07:34)
$ws.Range("A1").Value = "Datos actualizados a 10 de Mayo de 2020 a las 07:34"

# Row 66 now shows Hungria with refreshed figures (was Oman)
$ws.Range("A66").Value = "Hungria"
$ws.Range("B66").Value = 3263
$ws.Range("C66").Value = 50
$ws.Range("D66").Value = 933
$ws.Range("E66").Value = 1917
$ws.Range("F66").Value = 50
$ws.Range("G66").Value = 8
$ws.Range("H66").Value = 413

# Row 67 now shows Oman with the previous Oman figures (shifted down)
$ws.Range("A67").Value = "Oman"
$ws.Range("B67").Value = 3224
$ws.Range("C67").Value = 0
$ws.Range("D67").Value = 1068
$ws.Range("E67").Value = 2139
$ws.Range("F67").Value = 17
$ws.Range("G67").Value = 0
$ws.Range("H67").Value = 17

# Row 74: Uzbekistan figures refreshed
$ws.Range("B74").Value = 2387
$ws.Range("C74").Value = 38
$ws.Range("D74").Value = 1846
$ws.Range("E74").Value = 531
$ws.Range("F74").Value = 8
$ws.Range("G74").Value = 0
$ws.Range("H74").Value = 10

# Row 79: Bulgaria figures refreshed
$ws.Range("B79").Value = 1955
$ws.Range("C79").Value = 34
$ws.Range("D79").Value = 444
$ws.Range("E79").Value = 1421
$ws.Range("F79").Value = 65
$ws.Range("G79").Value = 0
$ws.Range("H79").Value = 90

# Row 95 now shows Kirguistan with refreshed figures (was Somalia)
$ws.Range("A95").Value = "Kirguistan"
$ws.Range("B95").Value = 1002
$ws.Range("C95").Value = 71
$ws.Range("D95").Value = 675
$ws.Range("E95").Value = 315
$ws.Range("F95").Value = 13
$ws.Range("G95").Value = 0
$ws.Range("H95").Value = 12

# Row 96 now shows Somalia with the previous Somalia figures (shifted down)
$ws.Range("A96").Value = "Somalia"
$ws.Range("B96").Value = 997
$ws.Range("C96").Value = 0
$ws.Range("D96").Value = 110
$ws.Range("E96").Value = 839
$ws.Range("F96").Value = 2
$ws.Range("G96").Value = 0
$ws.Range("H96").Value = 48

# Row 97 now shows Mayotte with the previous Mayotte figures (shifted down)
$ws.Range("A97").Value = "Mayotte"
$ws.Range("B97").Value = 988
$ws.Range("C97").Value = 0
$ws.Range("D97").Value = 492
$ws.Range("E97").Value = 485
$ws.Range("F97").Value = 7
$ws.Range("G97").Value = 0
$ws.Range("H97").Value = 11

# Row 98 now shows Guatemala with the previous Guatemala figures (shifted down)
$ws.Range("A98").Value = "Guatemala"
$ws.Range("B98").Value = 967
$ws.Range("C98").Value = 67
$ws.Range("D98").Value = 104
$ws.Range("E98").Value = 839
$ws.Range("F98").Value = 5
$ws.Range("G98").Value = 0
$ws.Range("H98").Value = 24

# Row 99 now shows Consejo Danes para los Refugiados with the previous figures (shifted down)
$ws.Range("A99").Value = "Consejo Danes para los Refugiados"
$ws.Range("B99").Value = 937
$ws.Range("C99").Value = 0
$ws.Range("D99").Value = 130
$ws.Range("E99").Value = 768
$ws.Range("F99").Value = 0
$ws.Range("G99").Value = 0
$ws.Range("H99").Value = 39

$wb.Save()
